# The "Data" worksheet repeats two header labels ("V1" and "V2") throughout
# the sheet (once per question block, in columns C and G respectively).
# The commit renames these labels to "KV1" and "KV2".
#
# Use a whole-cell (exact match) Find & Replace scoped to the "Data" sheet so
# that only cells whose entire content is exactly "V1" / "V2" are touched -
# this avoids accidentally matching unrelated strings that merely contain
# "V1"/"V2" as a substring (e.g. "V18_A1", "V19", ...).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# xlWhole = 2 -> match the entire cell contents only (no partial match)
# xlByRows = 1 -> search order (irrelevant for Replace, kept for clarity)
$ws.Cells.Replace("V1", "KV1", 2, 1, $false)
$ws.Cells.Replace("V2", "KV2", 2, 1, $false)
